# Added scale and lines connecting gpx points
# Populate the three previously-blank rows (161-163) on the "Translation"
# sheet with new text entries for the map scale indicator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 161: combined "<value><value>" text used to render the scale label
$ws.Cells.Item(161, 2).Value = "SingleUseId220"
$ws.Cells.Item(161, 3).Value = "Small"
$ws.Cells.Item(161, 4).Value = "Center"
$ws.Cells.Item(161, 5).Value = "LTR"
$ws.Cells.Item(161, 6).Value = "<value><value>"

# Row 162: unit suffix ("m" for meters) used by the scale bar
$ws.Cells.Item(162, 2).Value = "SingleUseId221"
$ws.Cells.Item(162, 3).Value = "Small"
$ws.Cells.Item(162, 4).Value = "Left"
$ws.Cells.Item(162, 5).Value = "LTR"
$ws.Cells.Item(162, 6).Value = "m"

# Row 163: default scale value (force text so "100" isn't stored as a number)
$ws.Cells.Item(163, 2).Value = "SingleUseId222"
$ws.Cells.Item(163, 3).Value = "Small"
$ws.Cells.Item(163, 4).Value = "Left"
$ws.Cells.Item(163, 5).Value = "LTR"
$ws.Cells.Item(163, 6).NumberFormat = "@"
$ws.Cells.Item(163, 6).Value = "100"
$ws.Cells.Item(163, 6).Style = "Normal"
